# "no db selected exception handled"
# Add a 4th worksheet (Sheet4) that mirrors the header row and the first
# data row of Sheet3, make it the active sheet/tab, and normalize the
# selection on every sheet to A1:D2 (keeping the active cell at A1).

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)
$sheet2 = $wb.Worksheets.Item(2)
$sheet3 = $wb.Worksheets.Item(3)

# Add the new sheet right after the last existing sheet (Sheet3), so it
# lands as "Sheet4" at the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet4 = $wb.Worksheets.Add($null, $lastSheet)

# Populate Sheet4 with the same header row as the other sheets plus the
# first data row of Sheet3.
$sheet4.Range("A1").Value = "first_name"
$sheet4.Range("B1").Value = "last_name"
$sheet4.Range("C1").Value = "age"
$sheet4.Range("D1").Value = "college"

$sheet4.Range("A2").Value = "sdjkasjd"
$sheet4.Range("B2").Value = "asS"
$sheet4.Range("C2").Value = "sSD"
$sheet4.Range("D2").Value = "sdasASD"

# Normalize the stored selection on every sheet to A1:D2 (active cell
# stays A1). Activating each sheet in turn updates tabSelected so only
# the last-activated sheet (Sheet4) ends up tabSelected, and becomes the
# workbook's active tab.
$sheet1.Activate()
$sheet1.Range("A1:D2").Select()

$sheet2.Activate()
$sheet2.Range("A1:D2").Select()

$sheet3.Activate()
$sheet3.Range("A1:D2").Select()

$sheet4.Activate()
$sheet4.Range("A1:D2").Select()
